$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.458.31"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "3.503.82"
$ws.Range("E3").Value = "  -3.35%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "605.84"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").Value = "149.11"
$ws.Range("E6").Value = "  -6.59%  "
$ws.Range("D7").Value = "3.504.88"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").Value = "7.58"
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("D12").Value = "0.428"
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  -5.86%  "
$ws.Range("D14").Value = "31.97"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "4.093.74"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").Value = "3.516.13"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "67.489.46"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "15.41"
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("D21").Value = "9.96"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "445.78"
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -3.56%  "
$ws.Range("D24").Value = "79.06"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "3.647.13"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "5.66"
$ws.Range("E27").Value = "  -4.67%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000122"
$ws.Range("E28").Value = "  -10.54%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "9.92"
$ws.Range("E29").Value = "  -7.90%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "8.50"
$ws.Range("E30").Value = "  -8.49%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.50"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.64"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.167"
$ws.Range("E33").Value = "  -7.01%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "25.62"
$ws.Range("E35").Value = "  -3.80%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "6.21"
$ws.Range("E36").Value = "  -6.34%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.498.09"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "7.98"
$ws.Range("E39").Value = "  -6.40%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.26"
$ws.Range("E41").Value = "  -7.59%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "174.16"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0896"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "5.41"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "30.93"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.899"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "46.73"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "7.58"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  -11.34%  "
